# Auto-generated Excel COM-interop script
# Applies updated profit-calculation values (columns H-N) for specific
# leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets,
# matching a scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3534.75
$ws.Range("I74").Value = 4023.4
$ws.Range("K74").Value = 4023.4
$ws.Range("M74").Value = -3087.4

$ws.Range("H77").Value = 3534.75
$ws.Range("I77").Value = 4023.4
$ws.Range("K77").Value = 20117
$ws.Range("M77").Value = -15437

$ws.Range("H82").Value = 4928
$ws.Range("I82").Value = 609.1429000000001
$ws.Range("K82").Value = 1827.4287
$ws.Range("M82").Value = -1421.4287

$ws.Range("H85").Value = 4928
$ws.Range("I85").Value = 609.1429000000001
$ws.Range("K85").Value = 1827.4287
$ws.Range("M85").Value = -423.4287000000002

$ws.Range("H86").Value = 100001950
$ws.Range("I86").Value = 111113170
$ws.Range("K86").Value = 111113170
$ws.Range("M86").Value = -111112047

$ws.Range("H89").Value = 100001950
$ws.Range("I89").Value = 111113170
$ws.Range("K89").Value = 555565850
$ws.Range("M89").Value = -555560234

$ws.Range("H100").Value = 528.8
$ws.Range("I100").Value = 471.93332
$ws.Range("J100").Value = 699.4
$ws.Range("K100").Value = 471.93332
$ws.Range("L100").Value = 699.4
$ws.Range("M100").Value = 69.06668000000002
$ws.Range("N100").Value = -1781.4

$ws.Range("H112").Value = 3005.4595
$ws.Range("J112").Value = 3075.0833
$ws.Range("L112").Value = 9225.249899999999
$ws.Range("N112").Value = -11441.2499

$ws.Range("H135").Value = 7553.385
$ws.Range("I135").Value = 1275.625
$ws.Range("K135").Value = 11480.625
$ws.Range("M135").Value = -8945.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12129.2
$ws.Range("I32").Value = 11539
$ws.Range("K32").Value = 11539
$ws.Range("M32").Value = -11252

$ws.Range("H45").Value = 83281
$ws.Range("I45").Value = 127216.664
$ws.Range("J45").Value = 4196.8
$ws.Range("K45").Value = 127216.664
$ws.Range("L45").Value = 4196.8
$ws.Range("M45").Value = -126839.664
$ws.Range("N45").Value = -4950.8

$ws.Range("H102").Value = 880.7931
$ws.Range("I102").Value = 786.95654
$ws.Range("K102").Value = 786.95654
$ws.Range("M102").Value = 835.04346

$ws.Range("H132").Value = 45975.855
$ws.Range("I132").Value = 84651.71000000001
$ws.Range("J132").Value = 7300
$ws.Range("K132").Value = 253955.13
$ws.Range("L132").Value = 21900
$ws.Range("M132").Value = -251425.13
$ws.Range("N132").Value = -26960

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5059.1
$ws.Range("I86").Value = 3984.3333
$ws.Range("K86").Value = 3984.3333
$ws.Range("M86").Value = -2861.3333

$ws.Range("H89").Value = 5059.1
$ws.Range("I89").Value = 3984.3333
$ws.Range("K89").Value = 19921.6665
$ws.Range("M89").Value = -14305.6665

$ws.Range("H105").Value = 4900.4443
$ws.Range("I105").Value = 3221.8
$ws.Range("K105").Value = 3221.8
$ws.Range("M105").Value = -1474.8

$ws.Range("H132").Value = 96648.53
$ws.Range("J132").Value = 96648.53
$ws.Range("L132").Value = 96648.53
$ws.Range("N132").Value = -106768.53

$ws.Range("H134").Value = 2974.468
$ws.Range("I134").Value = 2902.8462
$ws.Range("J134").Value = 3063.1428
$ws.Range("K134").Value = 8708.5386
$ws.Range("L134").Value = 9189.428400000001
$ws.Range("M134").Value = -6173.5386
$ws.Range("N134").Value = -14259.4284

$ws.Range("H140").Value = 246249.75
$ws.Range("J140").Value = 300000
$ws.Range("L140").Value = 300000
$ws.Range("N140").Value = -310360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2636.5
$ws.Range("I22").Value = 2636.5
$ws.Range("K22").Value = 2636.5
$ws.Range("M22").Value = -2286.5

$ws.Range("H31").Value = 18185238
$ws.Range("I31").Value = 20411070
$ws.Range("K31").Value = 20411070
$ws.Range("M31").Value = -20410775

$ws.Range("H34").Value = 18185238
$ws.Range("I34").Value = 20411070
$ws.Range("K34").Value = 20411070
$ws.Range("M34").Value = -20410868

$ws.Range("H62").Value = 42386.1
$ws.Range("I62").Value = 20555.285
$ws.Range("K62").Value = 20555.285
$ws.Range("M62").Value = -19931.285

$ws.Range("H65").Value = 42386.1
$ws.Range("I65").Value = 20555.285
$ws.Range("K65").Value = 102776.425
$ws.Range("M65").Value = -99656.425

$ws.Range("H134").Value = 1802.1212
$ws.Range("I134").Value = 1231.3572
$ws.Range("K134").Value = 3694.0716
$ws.Range("M134").Value = -1159.0716

$ws.Range("H141").Value = 154373.7
$ws.Range("J141").Value = 178388
$ws.Range("L141").Value = 178388
$ws.Range("N141").Value = -188748

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 8548901
$ws.Range("I131").Value = 822
$ws.Range("J131").Value = 10754857
$ws.Range("K131").Value = 2466
$ws.Range("L131").Value = 32264571
$ws.Range("M131").Value = 2574
$ws.Range("N131").Value = -32274651

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 11113562
$ws.Range("I102").Value = 13515854
$ws.Range("K102").Value = 13515854
$ws.Range("M102").Value = -13514232

$ws.Range("H107").Value = 1058.1818
$ws.Range("J107").Value = 1374.75
$ws.Range("L107").Value = 1374.75
$ws.Range("N107").Value = -5214.75

$ws.Range("H122").Value = 282332.06
$ws.Range("I122").Value = 529280.9399999999
$ws.Range("J122").Value = 6330.353
$ws.Range("K122").Value = 1587842.82
$ws.Range("L122").Value = 18991.059
$ws.Range("M122").Value = -1585392.82
$ws.Range("N122").Value = -23891.059

$ws.Range("H126").Value = 3712.158
$ws.Range("I126").Value = 3778.7222
$ws.Range("K126").Value = 11336.1666
$ws.Range("M126").Value = -8866.1666

$ws.Range("H132").Value = 136323.06
$ws.Range("I132").Value = 169405.5
$ws.Range("J132").Value = 3993.3333
$ws.Range("K132").Value = 508216.5
$ws.Range("L132").Value = 11979.9999
$ws.Range("M132").Value = -505686.5
$ws.Range("N132").Value = -17039.9999

$ws.Range("H141").Value = 37464.25
$ws.Range("J141").Value = 37464.25
$ws.Range("L141").Value = 37464.25
$ws.Range("N141").Value = -47824.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5954.706
$ws.Range("I7").Value = 5397.2856
$ws.Range("J7").Value = 6855.154
$ws.Range("K7").Value = 5397.2856
$ws.Range("L7").Value = 6855.154
$ws.Range("M7").Value = -5285.2856
$ws.Range("N7").Value = -7079.154

$ws.Range("H68").Value = 4398.2
$ws.Range("I68").Value = 2997.3333
$ws.Range("K68").Value = 2997.3333
$ws.Range("M68").Value = -2248.3333

$ws.Range("H71").Value = 4398.2
$ws.Range("I71").Value = 2997.3333
$ws.Range("K71").Value = 14986.6665
$ws.Range("M71").Value = -11242.6665

$ws.Range("H126").Value = 5954.706
$ws.Range("I126").Value = 5397.2856
$ws.Range("J126").Value = 6855.154
$ws.Range("K126").Value = 16191.8568
$ws.Range("L126").Value = 20565.462
$ws.Range("M126").Value = -13721.8568
$ws.Range("N126").Value = -25505.462

$ws.Range("H136").Value = 5111.1665
$ws.Range("I136").Value = 2945.1904
$ws.Range("J136").Value = 6489.515
$ws.Range("K136").Value = 8835.5712
$ws.Range("L136").Value = 19468.545
$ws.Range("M136").Value = -6285.5712
$ws.Range("N136").Value = -24568.545

$ws.Range("H139").Value = 88715
$ws.Range("J139").Value = 88715
$ws.Range("L139").Value = 88715
$ws.Range("N139").Value = -98995

$ws.Range("H141").Value = 112857.5
$ws.Range("J141").Value = 112857.5
$ws.Range("L141").Value = 112857.5
$ws.Range("N141").Value = -123217.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10241.538
$ws.Range("I81").Value = 8178
$ws.Range("J81").Value = 12649
$ws.Range("K81").Value = 16356
$ws.Range("L81").Value = 25298
$ws.Range("M81").Value = -15295
$ws.Range("N81").Value = -27420

$ws.Range("H84").Value = 10241.538
$ws.Range("I84").Value = 8178
$ws.Range("J84").Value = 12649
$ws.Range("K84").Value = 81780
$ws.Range("L84").Value = 126490
$ws.Range("M84").Value = -76476
$ws.Range("N84").Value = -137098

$ws.Range("H138").Value = 98400
$ws.Range("I138").Value = 98400
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 98400
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -93260
$ws.Range("N138").ClearContents()

$ws.Range("H140").Value = 97209.5
$ws.Range("J140").Value = 97209.5
$ws.Range("L140").Value = 97209.5
$ws.Range("N140").Value = -107569.5
